# Auto-generated edit script for cs_upload_template.xlsx
# Incorporated UG into the tool: remove 5 discontinued Uganda products
# (L.O.P. FOLLOWING MACHINERY B/DOWN, RETENTION BOND, AVIATION, GOLFERS INSURANCE, LIVESTOCK)
# from the valid_product_names!C (Uganda) list, and refresh the customer_prod demo rows.

$wb = $excel.ActiveWorkbook
$wsValid = $wb.Worksheets.Item("valid_product_names")
$wsProd  = $wb.Worksheets.Item("customer_prod")

# New Uganda product list (column C), rows 2-97 -- 5 discontinued products removed,
# remaining products shifted up.
$ugandaProducts = @(
    'CANCER TREATMENT',
    'CHRONIC / HIV-AIDS / PRE-EXISTING DISEASE',
    'CONGENITAL DISEASE',
    'DENTAL',
    'DENTAL SURGERY',
    'FUNERAL',
    'GYNAECOLOGICAL SURGERY',
    'IMMUNIZATION',
    'IN-PATIENT',
    'MAXILLOFACIAL AND ORAL SURGERY',
    'OPTHALMOLOGY',
    'ORGAN TRANSPLANT',
    'OUT-PATIENT',
    'POST-HOSPITALIZATION',
    'PROSTHESIS APPLIANCES',
    'PSYCHIATRIC / PHYSIOTHERAPY',
    'RECONSTRUCTIVE SURGERY',
    'ROOM CHARGE DAY LIMIT',
    'CHRONIC / PRE-EXISTING / HIV / CANCER',
    'DENTAL BRIDGES',
    'DENTAL CROWNS',
    'DENTAL DENTURES',
    'ILLNESS',
    'LODGER`S FEE',
    'MATERNITY',
    'MATERNITY / CAESAREAN',
    'OPTICAL',
    'PSYCHIATRIC TREATMENT',
    'VISIT FEE',
    'ANNUAL MEDICAL CHECK-UP',
    'CHRONIC / PRE-EXISTING / HIV / MATERNITY',
    'CONGENITAL AND NEONATAL CONDITION',
    'MATERNITY / EMERGENCY CAESAREAN',
    'NEONATAL',
    'MATERNITY (CAESAREAN)',
    'SPECTACLE FRAME',
    'HIV / AIDS',
    'ROOM CHARGE',
    'CHRONIC / CONGENITAL',
    'ANNUAL MEDICAL CHECK-UP FOR MAIN MEMBER / EMPLOYEE',
    'EXCESS OF LOSS',
    'ANNUAL MEDICAL CHECK-UP FOR DEPANDANT',
    'CONTACT LENSES',
    'DENTAL CLEANING',
    'DENTAL EXTRACTIONS',
    'DENTAL FILLING',
    'DENTAL POLISHING',
    'DENTAL ROOT CANAL',
    'DENTAL SCALING',
    'EYE TEST',
    'VACCINATION',
    'MATERNITY RELATED AILMENT',
    'PHYSIOTHERAPY',
    'PERSONAL ACCIDENT',
    'ANTENATAL AND POSTNATAL COVER',
    'MEDICAL AIDS',
    'MATERNITY / MATERNITY RELATED AILMENT',
    'CHRONIC / PRE-EXISTING / CONGENITAL',
    'MOTOR COMMERCIAL',
    'BURGLARY',
    'FIRE COMMERCIAL',
    'MOTOR PRIVATE',
    'ALL RISKS',
    'ELECTRONIC EQUIPMENT',
    'MOTOR - COMESA (YELLOW CARDS)',
    'MOTOR THIRD PARTY',
    'EMPLOYER''S LIABILITY',
    'FIDELITY GUARANTEE',
    'PUBLIC LIABILITY',
    'CASH IN TRANSIT',
    'BID BOND',
    'MOTOR CYCLE',
    'PERFORMANCE BOND',
    'PROFESSIONAL LIABILITY',
    'ADVANCE PAYMENT BOND',
    'MACHINERY BREAKDOWN',
    'PLATE GLASS',
    'GOODS IN TRANSIT',
    'HOUSEHOLDERS (HHC)',
    'UAP ASSETS ALL RISKS',
    'CONTRACTORS PLANT AND MACHINERY',
    'CONTRACTORS ALL RISKS',
    'MARINE CARGO',
    'INDUSTRIAL ALL RISKS',
    'HOUSEOWNERS (HOC)',
    'BONDS (ALL TYPES)',
    'CUSTOMS BOND',
    'POLISURE',
    'MARINE OPEN COVER',
    'TRAVELLERS INSURANCE',
    'FIRE CONSEQUENTIAL LOSSES',
    'ERECTION ALL RISKS',
    'SCHOOL''S COMPREHENSIVE',
    'MOTOR TRADE',
    'SURETY UNDERTAKING',
    'MARINE HULL'
)

$row = 2
foreach ($product in $ugandaProducts) {
    $wsValid.Cells.Item($row, 3).Value = $product
    $row = $row + 1
}

# The list shrank from 101 to 96 entries (C2:C102 -> C2:C97); clear the now-unused tail.
$wsValid.Range("C98:C102").ClearContents() | Out-Null

# Update the named range "Uganda" to match the new, shorter list.
$wb.Names.Item("Uganda").RefersTo = "=valid_product_names!`$C`$2:`$C`$97"

# Refresh the customer_prod sample/demo rows (A2:A6) to reflect the updated product set.
$wsProd.Range("A2").Value = 'CONGENITAL DISEASE'
$wsProd.Range("A3").Value = 'GYNAECOLOGICAL SURGERY'
$wsProd.Range("A4").Value = 'DENTAL SURGERY'
$wsProd.Range("A5").ClearContents() | Out-Null
$wsProd.Range("A6").ClearContents() | Out-Null

# Move the active selections to match the edited workbook state.
$wsProd.Range("A3").Select() | Out-Null
$wsValid.Range("C1").Select() | Out-Null
